# "Generate Report for Handoff"
# Flip the three localization-status worksheets from "handed back" state to
# "ready for handoff" state: update the status text + refresh the relevant
# timestamp, then let the status columns re-fit to the (shorter) new text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

# New column width (character units) the status columns settle on once the
# shorter "Ready for handoff" text is autosized.
$newColWidth = 16.33

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-09-04 15:03:16"
$wsOverview.Range("E1:F2").ColumnWidth = $newColWidth

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-09-04 15:03:12"
$wsZhCn.Range("C1:C2").ColumnWidth = $newColWidth

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-09-04 15:03:16"
$wsDeDe.Range("C1:C2").ColumnWidth = $newColWidth
